$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected, so it must be unprotected before any cell can
# be written, then re-protected afterwards to restore its original state.
$ws.Unprotect()

# Bump the "as of" date in the confidentiality / disclosure footnote.
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Refresh the model holdings weights / percent-change figures.
$ws.Range("D2").Value = 0.8470313387914886
$ws.Range("E2").Value = 0.00186358553857624
$ws.Range("D3").Value = 0.1529686612085114
$ws.Range("E3").Value = 0.0001834525775088469
$ws.Range("E4").Value = 0.001606577848869417

$ws.Protect()
